# Adds an "admin" category column to the project lookup table (Table1 on
# Sheet1): new table column "admin" with per-project values of
# Outsourced / In-house / Unknown, then re-applies the "phase" filter
# (Onboarded, Original) so the Phase 3 rows stay hidden behind the filter.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# 1. Add the new "admin" column to the table (expands Table1 / autofilter
#    ref from A1:N46 to A1:O46, bumps tableColumns count to 15).
$newCol = $lo.ListColumns.Add()
$lo.HeaderRowRange.Item(15).Value = "admin"

# 2. Populate the admin classification for each "Original"/"Onboarded"
#    project (rows 2-21). Phase 3 projects (rows 22-46) are left blank.
$adminValues = @{
    2  = "Outsourced"
    3  = "Unknown"
    4  = "Outsourced"
    5  = "In-house"
    6  = "In-house"
    7  = "Outsourced"
    8  = "Unknown"
    9  = "Outsourced"
    10 = "Unknown"
    11 = "Outsourced"
    12 = "In-house"
    13 = "Outsourced"
    14 = "In-house"
    15 = "Outsourced"
    16 = "Outsourced"
    17 = "Outsourced"
    18 = "Outsourced"
    19 = "In-house"
    20 = "In-house"
    21 = "Outsourced"
}
foreach ($r in $adminValues.Keys) {
    $ws.Cells.Item($r, 15).Value = $adminValues[$r]
}

# 3. Drop the cached "bestFit" auto-width flag on the columns that already
#    had an explicit width (C:N) - same widths, no longer flagged bestFit.
for ($c = 3; $c -le 14; $c++) {
    $col = $ws.Columns.Item($c)
    $col.ColumnWidth = $col.ColumnWidth
}
# Give the new admin column (O) a best-fit width of its own.
$ws.Columns.Item(15).AutoFit()

# 4. Re-apply the phase autofilter to Table1 (column 3 = "phase") so only
#    "Onboarded"/"Original" rows show - this hides rows 22:46 (Phase 3).
$lo.Range.AutoFilter(3, @("Onboarded", "Original"), 7)

# 5. Leave the active selection on B2, matching the saved view state.
$ws.Range("B2").Select()
